$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Invalid): rows 3 and 4 become 1
$ws.Range("G3:G4").Value = 1

# Column H (Absent): rows 3 through 18 become 1
$ws.Range("H3:H18").Value = 1
